$wb = $excel.ActiveWorkbook

# Helper: force a cell to hold literal TEXT (never auto-coerced to a
# number) without touching NumberFormat / styles. We write a formula
# that evaluates to the literal string, then convert that formula to
# its value in place (copy + paste-values). The end result is a plain
# text cell (shared string), no residual formula, no new cell style.
function Set-TextValue($range, $text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q1" sheet right before "总计" (it becomes the
#    4th sheet, pushing "总计" to the 5th / last position).
# ---------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalWs)
$q1.Name = "2022-Q1"

# Copy header formatting (style s=2: bold + bordered) from an existing
# quarter sheet so the new header row matches the others' look.
$srcQ = $wb.Worksheets.Item("2021-Q4")
$srcQ.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Copy the index-column style (s=2) used on the source sheet's A2:A3.
$srcQ.Range("A2:A3").Copy()
$q1.Range("A2:A3").PasteSpecial(-4122)

$q1.Range("A2").Value = 0
Set-TextValue $q1.Range("B2") "006143"
$q1.Range("C2").Value = "恒生前海中证质量成长低波动指数A"
Set-TextValue $q1.Range("D2") "0.06"
Set-TextValue $q1.Range("E2") "94.34"
Set-TextValue $q1.Range("F2") "3.02"
Set-TextValue $q1.Range("G2") "0.0018"
$q1.Range("H2").Value = 5

$q1.Range("A3").Value = 1
Set-TextValue $q1.Range("B3") "006144"
$q1.Range("C3").Value = "恒生前海中证质量成长低波动指数C"
Set-TextValue $q1.Range("D3") "0.01"
Set-TextValue $q1.Range("E3") "94.34"
Set-TextValue $q1.Range("F3") "3.02"
Set-TextValue $q1.Range("G3") "0.0003"
$q1.Range("H3").Value = 5

# ---------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: add a new top data row for
#    2022-Q1, shifting the existing quarter rows down by one and
#    renumbering the index column (A).
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

$tot.Range("A4").Copy()
$tot.Range("A5").PasteSpecial(-4122)
$tot.Range("A5").Value = 3
$tot.Range("B5").Value = $tot.Range("B4").Value2
$tot.Range("C5").Value = $tot.Range("C4").Value2
$tot.Range("D5").Value = $tot.Range("D4").Value2

$tot.Range("A4").Value = 2
$tot.Range("B4").Value = $tot.Range("B3").Value2
$tot.Range("C4").Value = $tot.Range("C3").Value2
$tot.Range("D4").Value = $tot.Range("D3").Value2

$tot.Range("A3").Value = 1
$tot.Range("B3").Value = $tot.Range("B2").Value2
$tot.Range("C3").Value = $tot.Range("C2").Value2
$tot.Range("D3").Value = $tot.Range("D2").Value2

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 2
$tot.Range("D2").Value = 0
